$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Row($r1, $r2) {
    $cols = @("D", "K", "L", "M", "N", "O", "P", "Q", "S")
    foreach ($col in $cols) {
        $cell1 = $ws.Range($col + $r1)
        $cell2 = $ws.Range($col + $r2)
        $v1 = $cell1.Value()
        $v2 = $cell2.Value()
        $cell1.Value = $v2
        $cell2.Value = $v1
    }
}

Swap-Row 2 4
Swap-Row 3 5
